# Update main_list.xlsx: replace the single data row (2002) with a new
# 2018 entry, and remove the other historical rows (2003-2017) so that
# only the header row and the new 2018 row remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove rows 3 through 17 (the old extra data rows), keeping header (row1)
# and the first data row (row2) which we will overwrite below.
$ws.Range("A3:B17").EntireRow.Delete()

# Set the remaining data row to the new 2018 entry.
$ws.Range("A2").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2018.xlsx"
$ws.Range("B2").Value = "2018"

# Update the active selection to match the target state.
$ws.Range("B2").Select()
